# Borrar puntos al final de frases en viñetas para versionas Full CV
#
# The "education" sheet lists, in column E, a set of bullet-point remarks
# attached to each degree entry. This removes the trailing full stop ('.')
# from each of those bullet sentences (the one that already had no trailing
# period, "Research project: 4.90/5.00", is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

$cells = @("E2", "E3", "E4", "E5", "E6", "E7")

foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $text = $cell.Value()
    if ($text.EndsWith(".")) {
        $cell.Value = $text.Substring(0, $text.Length - 1)
    }
}

$ws.Range("E7").Select()
